# Update automàtic: dades i banners [2026-02-20 09:49]
# Refresh the DATA_EXTRACCIO (column H) timestamps on the "Dades_Període" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Període")

$ws.Range("H2").Value = "2026-02-20 09:49:16"
$ws.Range("H3").Value = "2026-02-20 09:49:18"
$ws.Range("H4").Value = "2026-02-20 09:49:18"
$ws.Range("H5").Value = "2026-02-20 09:49:18"
$ws.Range("H6").Value = "2026-02-20 09:49:18"
